# Add new logbook entry "LiCl_018" as row 20 on Sheet1.
# Row 20 previously existed only as an empty placeholder row; this fills
# it in with the measurement data for run LiCl_018 (the data that sits
# between LiCl_017 / row 19 and LiCl_019 / row 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B20").Value = 734
$ws.Range("C20").Value = 733.5
$ws.Range("D20").Value = 1477.7355500000001
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 100
$ws.Range("G20").Value = 1757.5

$ws.Range("I20").Value = 3550
$ws.Range("J20").Value = 1162
$ws.Range("K20").Value = 0.034815526744208854
$ws.Range("L20").Value = 101966
$ws.Range("M20").Value = 1748.6590000000001

$ws.Range("P20").Value = 3214
$ws.Range("Q20").Value = 932
$ws.Range("R20").Value = 0.078749418077573324
$ws.Range("S20").Value = 40813
$ws.Range("T20").Value = 1752.3420000000001

$ws.Range("W20").Value = 69.5
$ws.Range("X20").Value = 90
$ws.Range("Y20").Value = "LiCl_018"

# Restore the current selection to match the saved workbook view.
$ws.Range("A1:AG59").Select()
